$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells text value, forcing Text number-format first when the
# string would otherwise be auto-parsed by Excel as a number (e.g. "1.001").
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Update Price (D) and Volume(1h) (E) columns for changed rows
$ws.Range("D2").Value = '30.768.55'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.911.27'
$ws.Range("E3").Value = '  +0.91%  '
Set-TextValue $ws.Range("D4") '1.001'
$ws.Range("E4").Value = '  -0.10%  '
Set-TextValue $ws.Range("D5") '239.92'
$ws.Range("E5").Value = '  -0.71%  '
Set-TextValue $ws.Range("D6") '1.001'
$ws.Range("E6").Value = '  +0.06%  '
Set-TextValue $ws.Range("D7") '0.4922'
$ws.Range("E7").Value = '  +0.05%  '
Set-TextValue $ws.Range("D8") '0.2960'
$ws.Range("E8").Value = '  +0.96%  '
Set-TextValue $ws.Range("D9") '0.06741'
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("D10").Value = '1.904.53'
$ws.Range("E10").Value = '  +0.55%  '
Set-TextValue $ws.Range("D11") '17.05'
$ws.Range("E11").Value = '  -0.59%  '
Set-TextValue $ws.Range("D12") '0.07362'
$ws.Range("E12").Value = '  +1.46%  '
Set-TextValue $ws.Range("D13") '5.166'
$ws.Range("E13").Value = '  +2.90%  '
Set-TextValue $ws.Range("D14") '88.30'
$ws.Range("E14").Value = '  -2.62%  '
Set-TextValue $ws.Range("D15") '0.6707'
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").Value = '30.719.54'
$ws.Range("E16").Value = '  +0.15%  '
Set-TextValue $ws.Range("D17") '0.000007911'
$ws.Range("E17").Value = '  -0.45%  '
Set-TextValue $ws.Range("D18") '13.49'
$ws.Range("E18").Value = '  +3.09%  '
Set-TextValue $ws.Range("D19") '1.001'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = '2.163.44'
$ws.Range("E20").Value = '  +1.05%  '
Set-TextValue $ws.Range("D21") '5.314'
$ws.Range("E21").Value = '  +10.82%  '
Set-TextValue $ws.Range("D22") '1.000'
$ws.Range("E22").Value = '  -0.14%  '
Set-TextValue $ws.Range("D23") '195.93'
$ws.Range("E23").Value = '  +3.41%  '
Set-TextValue $ws.Range("D24") '6.253'
$ws.Range("E24").Value = '  +3.04%  '
Set-TextValue $ws.Range("D25") '9.631'
$ws.Range("E25").Value = '  +3.02%  '
Set-TextValue $ws.Range("D26") '162.93'
$ws.Range("E26").Value = '  +3.95%  '
Set-TextValue $ws.Range("D27") '18.59'
$ws.Range("E27").Value = '  -0.82%  '
Set-TextValue $ws.Range("D28") '1.946'
$ws.Range("E28").Value = '  +2.99%  '
Set-TextValue $ws.Range("D29") '1.468'
$ws.Range("E29").Value = '  +4.83%  '
Set-TextValue $ws.Range("D30") '4.406'
$ws.Range("E30").Value = '  +3.09%  '
Set-TextValue $ws.Range("D31") '0.09124'
$ws.Range("E31").Value = '  +0.47%  '
Set-TextValue $ws.Range("D32") '4.055'
$ws.Range("E32").Value = '  +1.71%  '
Set-TextValue $ws.Range("D33") '0.05242'
$ws.Range("E33").Value = '  +0.27%  '
Set-TextValue $ws.Range("D34") '0.7396'
$ws.Range("E34").Value = '  +0.43%  '
Set-TextValue $ws.Range("D35") '1.112'
$ws.Range("E35").Value = '  +0.77%  '
Set-TextValue $ws.Range("D36") '2.723'
$ws.Range("E36").Value = '  -1.21%  '
Set-TextValue $ws.Range("D38") '2.712'
$ws.Range("E38").Value = '  +1.23%  '
Set-TextValue $ws.Range("D39") '0.9212'
$ws.Range("E39").Value = '  -1.03%  '
Set-TextValue $ws.Range("D40") '2.076'
$ws.Range("E40").Value = '  -2.08%  '
Set-TextValue $ws.Range("D41") '74.43'
$ws.Range("E41").Value = '  +29.25%  '
Set-TextValue $ws.Range("D42") '0.4441'
$ws.Range("E42").Value = '  +1.32%  '
Set-TextValue $ws.Range("D43") '106.67'
$ws.Range("E43").Value = '  +1.64%  '
Set-TextValue $ws.Range("D44") '5.918'
$ws.Range("E44").Value = '  +3.40%  '
Set-TextValue $ws.Range("D46") '0.1387'
$ws.Range("E46").Value = '  +2.93%  '
Set-TextValue $ws.Range("D47") '7.591'
$ws.Range("E47").Value = '  +1.02%  '
Set-TextValue $ws.Range("D50") '0.05850'
$ws.Range("E50").Value = '  -0.13%  '
Set-TextValue $ws.Range("D51") '0.3995'
$ws.Range("E51").Value = '  +1.74%  '

# Row 37 (VeChain): only Price changes, Volume(1h) stays the same
Set-TextValue $ws.Range("D37") '0.01820'

# Row 45 (PaxDollar): only Volume(1h) changes, Price stays the same
$ws.Range("E45").Value = '  -0.07%  '

# Rows 48 and 49 swap coin ordering (Elrond <-> EnergySwap) with refreshed data
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D48") '9.092'
$ws.Range("E48").Value = '  +4.21%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range("D49") '35.34'
$ws.Range("E49").Value = '  +5.17%  '
